# Apply the new built-in table style to the three financial tables
# (Balance sheet / Income statement slides) that currently use the
# custom "Table_0" style defined in tableStyles.xml.
#
# Table.Style is read-only in the PowerPoint object model - the style
# has to be changed with Table.ApplyStyle(styleId).

$p = $ppt.ActivePresentation

$oldStyleId = "{48D96741-1E92-4999-B031-7404B86E98C9}"
$newStyleId = "{3A388B44-C8D4-4BF4-8DB7-5AE19595C054}"

$targetSlideIndexes = 14, 15, 16

foreach ($slideIndex in $targetSlideIndexes) {
    $slide = $p.Slides.Item($slideIndex)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)

        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
